$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Name", "Fact", "Habitat", "Lifespan", "Diet"),
    @("Lion", "A lion's roar can be heard from 8 kilometers away.", "Grasslands, savannas", "10-14 years", "Carnivore"),
    @("Octopus", "An octopus has three hearts.", "Oceans, coral reefs", "1-2 years", "Carnivore"),
    @("Elephant", "Elephants can recognize themselves in a mirror.", "Forests, savannas", "60-70 years", "Herbivore"),
    @("Penguin", "Penguins can stay underwater for up to 27 minutes.", "Antarctic coasts", "15-20 years", "Carnivore"),
    @("Giraffe", "Giraffes have the same number of neck bones as humans.", "African savannas", "25 years", "Herbivore"),
    @("Kangaroo", "Kangaroos can leap up to 3 meters high.", "Australian bushland", "20 years", "Herbivore"),
    @("Sloth", "Sloths can hold their breath for up to 40 minutes.", "Central and South American rainforests", "20-30 years", "Herbivore"),
    @("Cheetah", "Cheetahs are the fastest land animals, reaching 112 km/h.", "African savannas", "10-12 years", "Carnivore"),
    @("Polar Bear", "Polar bears have black skin under their white fur.", "Arctic regions", "25-30 years", "Carnivore"),
    @("Chimpanzee", "Chimpanzees can use tools and solve problems.", "African forests", "40-50 years", "Omnivore")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
